{"js": "// Remove the trailing \"empty paragraph\" + \"Ver no Jupiter...\" paragraph +\n// \"\u00a9 2020 ... Creative Commons Attribution\" paragraph that followed the\n// SERAFINI bibliography entry.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\n// Find the index of the \"Ver no Jupiter...\" paragraph; the empty paragraph\n// immediately preceding it (which separates it from the SERAFINI reference)\n// is removed along with it and the following copyright paragraph.\nlet verIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === targetTexts[0]) {\n    verIdx = i;\n    break;\n  }\n}\n\nif (verIdx > -1) {\n  const toDelete = [];\n  // The blank paragraph right before \"Ver no Jupiter...\" (if present).\n  if (verIdx - 1 >= 0 && items[verIdx - 1].text === \"\") {\n    toDelete.push(items[verIdx - 1]);\n  }\n  toDelete.push(items[verIdx]);\n  if (verIdx + 1 < items.length && items[verIdx + 1].text === targetTexts[1]) {\n    toDelete.push(items[verIdx + 1]);\n  }\n  for (const p of toDelete) {\n    p.delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"empty paragraph\" + \"Ver no Jupiter...\" paragraph +\n# \"\u00a9 2020 ... Creative Commons Attribution\" paragraph that followed the\n# SERAFINI bibliography entry.\n$d = $word.ActiveDocument\n\n$verText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$verIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($text -eq $verText) {\n        $verIndex = $i\n        break\n    }\n}\n\nif ($verIndex -gt 0) {\n    $deleteIndices = New-Object System.Collections.ArrayList\n\n    # The paragraph right after \"Ver no Jupiter...\" \u2014 only delete it if it\n    # really is the copyright paragraph.\n    if ($verIndex + 1 -le $d.Paragraphs.Count) {\n        $nextText = $d.Paragraphs.Item($verIndex + 1).Range.Text.TrimEnd([char]13, [char]7)\n        if ($nextText -eq $copyrightText) {\n            [void]$deleteIndices.Add($verIndex + 1)\n        }\n    }\n\n    [void]$deleteIndices.Add($verIndex)\n\n    # The blank paragraph right before \"Ver no Jupiter...\" (if present).\n    if ($verIndex - 1 -ge 1) {\n        $prevText = $d.Paragraphs.Item($verIndex - 1).Range.Text.TrimEnd([char]13, [char]7)\n        if ($prevText -eq \"\") {\n            [void]$deleteIndices.Add($verIndex - 1)\n        }\n    }\n\n    # Delete from the highest index down so earlier indices stay valid.\n    $sorted = $deleteIndices | Sort-Object -Descending\n    foreach ($idx in $sorted) {\n        $d.Paragraphs.Item($idx).Range.Delete()\n    }\n}\n"}
